$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded. Insert a row above row 411
# (shifting the existing rows 411:483 down to 412:484) and populate it with
# the new data point, matching how the rest of the "Femacal de La Calera -
# Zanahoria" table is laid out.
$ws.Rows.Item(411).Insert()

$ws.Cells.Item(411, 1).Value = 3
$ws.Cells.Item(411, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(411, 3).Value = "Coquimbo"
$ws.Cells.Item(411, 4).Value = 44951
$ws.Cells.Item(411, 5).Value = 5
$ws.Cells.Item(411, 6).Value = 100114013
$ws.Cells.Item(411, 7).Value = "Zanahoria"
$ws.Cells.Item(411, 8).Value = "Sin especificar"
$ws.Cells.Item(411, 9).Value = "Primera"
$ws.Cells.Item(411, 10).Value = 600
$ws.Cells.Item(411, 11).Value = 11000
$ws.Cells.Item(411, 12).Value = 12000
$ws.Cells.Item(411, 13).Value = 11583
$ws.Cells.Item(411, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(411, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(411, 16).Value = 579
$ws.Cells.Item(411, 17).Value = 20
$ws.Cells.Item(411, 18).Value = "Hortaliza"
